$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AddressesIdea2")

# New helper column M: screen buffer size (320x256 pixels, 1bpp packed) and
# a label in column N identifying the row.
$ws.Range("M11").Formula = "=(320*256)/8"
$ws.Range("N11").Value = "Screen Size"

# Bank index selector values for the MMU bugfix (column N, rows 12-15)
$ws.Range("N12").Value = 0
$ws.Range("N13").Value = 1
$ws.Range("N14").Value = 2
$ws.Range("N15").Value = 3

# Column L: hex representation of the bank start address
$ws.Range("L12").Formula = "=DEC2HEX(M12)"
$ws.Range("L13:L15").Formula = "=DEC2HEX(M13)"

# Column M: computed bank start address based on the bank index in N
$ws.Range("M12").Formula = "=`$A`$12+(`$M`$11*N12)"
$ws.Range("M13:M15").Formula = "=`$A`$12+(`$M`$11*N13)"

# Leave the selection where the author ended up after the edit
$ws.Range("N21").Select()
